# Updated cryptos list on Sat Feb 18 15:51:11 UTC 2023 with GitHub Actions
#
# This script re-applies the refreshed price / 24h-volume figures scraped
# from coinranking.com, including the two pairs of rows (35/36 and 45/46)
# whose coins swapped rank position. Every target cell in the source sheet
# is a text cell (t="inlineStr"), so values are written back as text -
# NumberFormat is temporarily forced to "@" (Text) before the write and the
# cell Style is reset to "Normal" afterwards so no stray numeric coercion
# or leftover cell-format changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '24.711.23' }
    @{ Cell = 'E2'; Value = '  +2.66%  ' }
    @{ Cell = 'D3'; Value = '1.697.15' }
    @{ Cell = 'E3'; Value = '  +1.64%  ' }
    @{ Cell = 'D4'; Value = '0.9988' }
    @{ Cell = 'E4'; Value = '  +0.13%  ' }
    @{ Cell = 'D5'; Value = '316.76' }
    @{ Cell = 'E5'; Value = '  +2.22%  ' }
    @{ Cell = 'D6'; Value = '0.9995' }
    @{ Cell = 'E6'; Value = '  +0.23%  ' }
    @{ Cell = 'D7'; Value = '0.3951' }
    @{ Cell = 'E7'; Value = '  +1.78%  ' }
    @{ Cell = 'D8'; Value = '0.4042' }
    @{ Cell = 'E8'; Value = '  +1.18%  ' }
    @{ Cell = 'E9'; Value = '  +2.92%  ' }
    @{ Cell = 'D10'; Value = '1.0000' }
    @{ Cell = 'D11'; Value = '51.47' }
    @{ Cell = 'E11'; Value = '  -3.30%  ' }
    @{ Cell = 'D12'; Value = '0.08797' }
    @{ Cell = 'E12'; Value = '  +1.40%  ' }
    @{ Cell = 'D13'; Value = '7.233' }
    @{ Cell = 'E13'; Value = '  +6.13%  ' }
    @{ Cell = 'D14'; Value = '23.49' }
    @{ Cell = 'E14'; Value = '  +3.25%  ' }
    @{ Cell = 'D15'; Value = '8.172' }
    @{ Cell = 'E15'; Value = '  +11.58%  ' }
    @{ Cell = 'D16'; Value = '0.00001320' }
    @{ Cell = 'E16'; Value = '  +0.85%  ' }
    @{ Cell = 'D17'; Value = '1.692.92' }
    @{ Cell = 'E17'; Value = '  +1.51%  ' }
    @{ Cell = 'D18'; Value = '99.96' }
    @{ Cell = 'E18'; Value = '  +0.81%  ' }
    @{ Cell = 'D19'; Value = '0.06999' }
    @{ Cell = 'E19'; Value = '  +1.38%  ' }
    @{ Cell = 'E20'; Value = '  +3.27%  ' }
    @{ Cell = 'D21'; Value = '7.092' }
    @{ Cell = 'E21'; Value = '  +7.17%  ' }
    @{ Cell = 'D22'; Value = '0.9992' }
    @{ Cell = 'D23'; Value = '14.38' }
    @{ Cell = 'E23'; Value = '  +3.92%  ' }
    @{ Cell = 'D24'; Value = '24.682.11' }
    @{ Cell = 'E24'; Value = '  +2.52%  ' }
    @{ Cell = 'D25'; Value = '3.141' }
    @{ Cell = 'E25'; Value = '  +3.45%  ' }
    @{ Cell = 'D26'; Value = '2.341' }
    @{ Cell = 'E26'; Value = '  +1.61%  ' }
    @{ Cell = 'D27'; Value = '22.90' }
    @{ Cell = 'E27'; Value = '  +5.05%  ' }
    @{ Cell = 'D28'; Value = '162.04' }
    @{ Cell = 'E28'; Value = '  +1.26%  ' }
    @{ Cell = 'D29'; Value = '137.53' }
    @{ Cell = 'E29'; Value = '  +5.41%  ' }
    @{ Cell = 'D30'; Value = '5.206' }
    @{ Cell = 'E30'; Value = '  +1.31%  ' }
    @{ Cell = 'D31'; Value = '7.510' }
    @{ Cell = 'E31'; Value = '  +4.03%  ' }
    @{ Cell = 'D32'; Value = '1.881.26' }
    @{ Cell = 'E32'; Value = '  +1.66%  ' }
    @{ Cell = 'D33'; Value = '1.084' }
    @{ Cell = 'E33'; Value = '  -1.86%  ' }
    @{ Cell = 'D34'; Value = '0.08634' }
    @{ Cell = 'E34'; Value = '  +0.14%  ' }
    @{ Cell = 'B35'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D35'; Value = '7.132' }
    @{ Cell = 'E35'; Value = '  +0.92%  ' }
    @{ Cell = 'B36'; Value = 'FraxShare' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D36'; Value = '11.62' }
    @{ Cell = 'E36'; Value = '  +6.02%  ' }
    @{ Cell = 'D37'; Value = '0.2757' }
    @{ Cell = 'E37'; Value = '  +3.97%  ' }
    @{ Cell = 'D38'; Value = '1.928' }
    @{ Cell = 'E38'; Value = '  +0.44%  ' }
    @{ Cell = 'D39'; Value = '14.51' }
    @{ Cell = 'E39'; Value = '  -0.24%  ' }
    @{ Cell = 'D40'; Value = '0.09154' }
    @{ Cell = 'E40'; Value = '  +3.60%  ' }
    @{ Cell = 'E41'; Value = '  +7.58%  ' }
    @{ Cell = 'D42'; Value = '1.479' }
    @{ Cell = 'E42'; Value = '  +1.97%  ' }
    @{ Cell = 'E43'; Value = '  +1.51%  ' }
    @{ Cell = 'D44'; Value = '2.658' }
    @{ Cell = 'E44'; Value = '  +9.92%  ' }
    @{ Cell = 'B45'; Value = 'EnergySwap' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D45'; Value = '15.82' }
    @{ Cell = 'E45'; Value = '  +4.12%  ' }
    @{ Cell = 'B46'; Value = 'Decentraland' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D46'; Value = '0.7188' }
    @{ Cell = 'E46'; Value = '  +1.35%  ' }
    @{ Cell = 'D47'; Value = '4.224' }
    @{ Cell = 'E47'; Value = '  +2.66%  ' }
    @{ Cell = 'E48'; Value = '  +0.21%  ' }
    @{ Cell = 'D49'; Value = '141.10' }
    @{ Cell = 'E49'; Value = '  +1.26%  ' }
    @{ Cell = 'D50'; Value = '1.331' }
    @{ Cell = 'E50'; Value = '  +9.20%  ' }
    @{ Cell = 'D51'; Value = '0.07995' }
    @{ Cell = 'E51'; Value = '  +1.90%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
